# Auto-generated edit script: updates market-price / leve-profit columns (H:N)
# across several sheets, matching the scheduled-runner data refresh described
# in the commit diff. Values come straight from the target OOXML diff.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H19").Value = 774.0769
$ws.Range("I19").Value = 690.6
$ws.Range("J19").Value = 826.25
$ws.Range("K19").Value = 690.6
$ws.Range("L19").Value = 826.25
$ws.Range("M19").Value = -515.6
$ws.Range("N19").Value = -1176.25

$ws.Range("H129").Value = 1236.7778
$ws.Range("I129").Value = 0
$ws.Range("J129").Value = 1236.7778
$ws.Range("K129").Value = 0
$ws.Range("L129").Value = 3710.3334
$ws.Range("M129").ClearContents()
$ws.Range("N129").Value = -13710.3334

$ws.Range("H137").Value = 1315.9756
$ws.Range("I137").Value = 1188.2858
$ws.Range("J137").Value = 1591
$ws.Range("K137").Value = 3564.8574
$ws.Range("L137").Value = 4773
$ws.Range("M137").Value = -1014.8574
$ws.Range("N137").Value = -9873

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 3707.7292
$ws.Range("I32").Value = 3021.0967
$ws.Range("K32").Value = 3021.0967
$ws.Range("M32").Value = -2734.0967

$ws.Range("H74").Value = 3071.077
$ws.Range("I74").Value = 3378.1777
$ws.Range("J74").Value = 1096.8572
$ws.Range("K74").Value = 3378.1777
$ws.Range("L74").Value = 1096.8572
$ws.Range("M74").Value = -2504.1777
$ws.Range("N74").Value = -2844.8572

$ws.Range("H77").Value = 3071.077
$ws.Range("I77").Value = 3378.1777
$ws.Range("J77").Value = 1096.8572
$ws.Range("K77").Value = 16890.8885
$ws.Range("L77").Value = 5484.286
$ws.Range("M77").Value = -12522.8885
$ws.Range("N77").Value = -14220.286

$ws.Range("H122").Value = 1759.4193
$ws.Range("I122").Value = 1455.0714
$ws.Range("J122").Value = 4600
$ws.Range("K122").Value = 4365.2142
$ws.Range("L122").Value = 13800
$ws.Range("M122").Value = -1915.2142
$ws.Range("N122").Value = -18700

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H4").Value = 1316463.4
$ws.Range("I4").Value = 1666750.5
$ws.Range("J4").Value = 1154792.2
$ws.Range("K4").Value = 1666750.5
$ws.Range("L4").Value = 1154792.2
$ws.Range("M4").Value = -1666638.5
$ws.Range("N4").Value = -1155016.2

$ws.Range("H5").Value = 915.1429000000001
$ws.Range("I5").Value = 99.5
$ws.Range("J5").Value = 2002.6666
$ws.Range("K5").Value = 99.5
$ws.Range("L5").Value = 2002.6666
$ws.Range("M5").Value = 12.5
$ws.Range("N5").Value = -2226.6666

$ws.Range("H99").Value = 4934.4473
$ws.Range("I99").Value = 1922.238
$ws.Range("J99").Value = 8655.412
$ws.Range("K99").Value = 1922.238
$ws.Range("L99").Value = 8655.412
$ws.Range("M99").Value = -424.2380000000001
$ws.Range("N99").Value = -11651.412

$ws.Range("H107").Value = 1017.3125
$ws.Range("I107").Value = 697.55554
$ws.Range("J107").Value = 1428.4286
$ws.Range("K107").Value = 697.55554
$ws.Range("L107").Value = 1428.4286
$ws.Range("M107").Value = 1222.44446
$ws.Range("N107").Value = -5268.4286

$ws.Range("H126").Value = 4934.4473
$ws.Range("I126").Value = 1922.238
$ws.Range("J126").Value = 8655.412
$ws.Range("K126").Value = 5766.714
$ws.Range("L126").Value = 25966.236
$ws.Range("M126").Value = -3296.714
$ws.Range("N126").Value = -30906.236

$ws.Range("H132").Value = 2858.2173
$ws.Range("I132").Value = 1345.1666
$ws.Range("J132").Value = 4508.8184
$ws.Range("K132").Value = 4035.4998
$ws.Range("L132").Value = 13526.4552
$ws.Range("M132").Value = -1505.4998
$ws.Range("N132").Value = -18586.4552

$ws.Range("H134").Value = 1897.3182
$ws.Range("I134").Value = 1685.6216
$ws.Range("J134").Value = 3016.2856
$ws.Range("K134").Value = 5056.864799999999
$ws.Range("L134").Value = 9048.856800000001
$ws.Range("M134").Value = -2521.864799999999
$ws.Range("N134").Value = -14118.8568

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H120").Value = 5108.3335
$ws.Range("I120").Value = 4515
$ws.Range("J120").Value = 5850
$ws.Range("K120").Value = 13545
$ws.Range("L120").Value = 17550
$ws.Range("M120").Value = -8707
$ws.Range("N120").Value = -27226

$ws.Range("H129").Value = 2009.64
$ws.Range("I129").Value = 741.6667
$ws.Range("J129").Value = 2722.875
$ws.Range("K129").Value = 2225.0001
$ws.Range("L129").Value = 8168.625
$ws.Range("M129").Value = 2774.9999
$ws.Range("N129").Value = -18168.625

$ws.Range("H131").Value = 6389.95
$ws.Range("I131").Value = 742.8570999999999
$ws.Range("J131").Value = 9430.691999999999
$ws.Range("K131").Value = 2228.5713
$ws.Range("L131").Value = 28292.076
$ws.Range("M131").Value = 2811.4287
$ws.Range("N131").Value = -38372.076

$ws.Range("H132").Value = 1817.0834
$ws.Range("I132").Value = 1410.4546
$ws.Range("J132").Value = 2161.1538
$ws.Range("K132").Value = 12694.0914
$ws.Range("L132").Value = 19450.3842
$ws.Range("M132").Value = -10164.0914
$ws.Range("N132").Value = -24510.3842

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H97").Value = 1258.95
$ws.Range("I97").Value = 1284.2142
$ws.Range("J97").Value = 1200
$ws.Range("K97").Value = 1284.2142
$ws.Range("L97").Value = 1200
$ws.Range("M97").Value = -788.2141999999999
$ws.Range("N97").Value = -2192

$ws.Range("H132").Value = 2627.125
$ws.Range("I132").Value = 2673.6924
$ws.Range("J132").Value = 2572.0908
$ws.Range("K132").Value = 8021.0772
$ws.Range("L132").Value = 7716.2724
$ws.Range("M132").Value = -5491.0772
$ws.Range("N132").Value = -12776.2724

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H95").Value = 0
$ws.Range("J95").Value = 0
$ws.Range("L95").Value = 0
$ws.Range("N95").ClearContents()

$ws.Range("H107").Value = 2652.4614
$ws.Range("I107").Value = 2349.2
$ws.Range("J107").Value = 3663.3333
$ws.Range("K107").Value = 7047.599999999999
$ws.Range("L107").Value = 10989.9999
$ws.Range("M107").Value = -5127.599999999999
$ws.Range("N107").Value = -14829.9999

$ws.Range("H132").Value = 1858.2157
$ws.Range("I132").Value = 1494.7954
$ws.Range("J132").Value = 4142.5713
$ws.Range("K132").Value = 4484.3862
$ws.Range("L132").Value = 12427.7139
$ws.Range("M132").Value = -1954.3862
$ws.Range("N132").Value = -17487.7139

$ws.Range("H136").Value = 5465603
$ws.Range("I136").Value = 8333974
$ws.Range("J136").Value = 2039.5238
$ws.Range("K136").Value = 25001922
$ws.Range("L136").Value = 6118.5714
$ws.Range("M136").Value = -24999372
$ws.Range("N136").Value = -11218.5714
